$d = $word.ActiveDocument

# Pull the whole package as OOXML text so we can make a precise, surgical
# run-level edit identical to the authored change (splitting the "8x8 BitMap"
# label into "8x8 Bi"+"tPic" for the first textbox pair and "8x8 Bit"+"Pic"
# for the second textbox pair, which both render as "8x8 BitPic").
$xml = $d.Content.WordOpenXML

$runProps = '<w:rPr><w:rFonts w:cs="Al Bayan"/><w:b/><w:bCs/></w:rPr>'
$pattern = '<w:r w:rsidRPr="009861D2">' + $runProps + '<w:t xml:space="preserve">8x8 </w:t></w:r><w:r w:rsidRPr="009861D2">' + $runProps + '<w:t>BitMap</w:t></w:r>'

$replacementA = '<w:r w:rsidRPr="009861D2">' + $runProps + '<w:t>8x8 Bi</w:t></w:r><w:r w:rsidRPr="009861D2">' + $runProps + '<w:t>tPic</w:t></w:r>'
$replacementB = '<w:r w:rsidRPr="009861D2">' + $runProps + '<w:t>8x8 Bit</w:t></w:r><w:r w:rsidRPr="009861D2">' + $runProps + '<w:t>Pic</w:t></w:r>'

$result = ""
$cursor = 0
$occurrence = 0

while ($true) {
    $idx = $xml.IndexOf($pattern, $cursor)
    if ($idx -lt 0) { break }

    $occurrence++
    if ($occurrence -le 2) {
        $replacement = $replacementA
    } else {
        $replacement = $replacementB
    }

    $result += $xml.Substring($cursor, $idx - $cursor)
    $result += $replacement
    $cursor = $idx + $pattern.Length
}
$result += $xml.Substring($cursor)

if ($occurrence -ne 4) {
    throw ("Expected 4 occurrences of the '8x8 '/'BitMap' run pair, found " + $occurrence)
}

$null = $d.Content.InsertXML($result)
Write-Output ("Replaced occurrences: " + $occurrence)
